# Auto-generated edit script: updates the "Price" (D) and "Volume(1h)" (E)
# columns of the crypto listing on Sheet1 to reflect refreshed market data.
#
# Some D-column values (e.g. "212.47", "0.0616") look like plain numbers.
# Assigning them straight to Range.Value would make Excel auto-convert them
# to numeric cells, which would not match the source data (these columns
# are plain text). To keep them as text we prefix such values with a
# leading apostrophe, exactly like a user typing '212.47 into a cell in
# the Excel UI; Excel stores the cell as text (quotePrefix) and the
# apostrophe itself is not part of the stored/displayed value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.322.80'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '1.623.02'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '''212.47'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  +1.29%  '
$ws.Range('E8').Value = '  +1.58%  '
$ws.Range('D9').Value = '''0.0616'
$ws.Range('E9').Value = '  +1.09%  '
$ws.Range('D10').Value = '''18.98'
$ws.Range('E10').Value = '  +6.19%  '
$ws.Range('D11').Value = '''0.0816'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '1.848.78'
$ws.Range('E12').Value = '  +1.87%  '
$ws.Range('D13').Value = '1.626.32'
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('D14').Value = '''4.02'
$ws.Range('E14').Value = '  +0.65%  '
$ws.Range('D15').Value = '''0.520'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '26.323.08'
$ws.Range('E16').Value = '  +1.44%  '
$ws.Range('D17').Value = '''62.63'
$ws.Range('E17').Value = '  +4.50%  '
$ws.Range('D18').Value = '0.0₃0730'
$ws.Range('E18').Value = '  +1.35%  '
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').Value = '''203.19'
$ws.Range('E20').Value = '  +2.07%  '
$ws.Range('D21').Value = '''4.31'
$ws.Range('E21').Value = '  +2.23%  '
$ws.Range('D22').Value = '''9.37'
$ws.Range('E22').Value = '  +1.68%  '
$ws.Range('D23').Value = '''6.06'
$ws.Range('E23').Value = '  +1.31%  '
$ws.Range('D24').Value = '''1.92'
$ws.Range('E24').Value = '  +7.03%  '
$ws.Range('D25').Value = '''143.76'
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.08%  '
$ws.Range('D28').Value = '''15.23'
$ws.Range('E28').Value = '  +1.24%  '
$ws.Range('E29').Value = '  +2.13%  '
$ws.Range('D30').Value = '''0.0529'
$ws.Range('E30').Value = '  +11.76%  '
$ws.Range('E31').Value = '  +0.90%  '
$ws.Range('E32').Value = '  +2.81%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('E34').Value = '  +2.25%  '
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('D36').Value = '1.180.44'
$ws.Range('E36').Value = '  +4.97%  '
$ws.Range('E37').Value = '  +1.65%  '
$ws.Range('D38').Value = '''0.812'
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('D41').Value = '''0.497'
$ws.Range('E41').Value = '  +2.04%  '
$ws.Range('D42').Value = '''0.790'
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('D43').Value = '''5.35'
$ws.Range('E43').Value = '  +5.42%  '
$ws.Range('D44').Value = '1.759.34'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = '''93.46'
$ws.Range('E45').Value = '  +1.18%  '
$ws.Range('D46').Value = '0.0₆0105'
$ws.Range('E46').Value = '  +14.54%  '
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('D48').Value = '''54.01'
$ws.Range('E48').Value = '  +1.49%  '
$ws.Range('E49').Value = '  +1.12%  '
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('E51').Value = '  -0.46%  '
